$wb = $excel.ActiveWorkbook

# Rename "ArrayTryCode" sheet to "TryCode"
$wsTryCode = $wb.Worksheets.Item("ArrayTryCode")
$wsTryCode.Name = "TryCode"

# Format B8 and B10 on PracticeQns as Text (numFmtId 49)
$wsPracticeQns = $wb.Worksheets.Item("PracticeQns")
$wsPracticeQns.Range("B8").NumberFormat = "@"
$wsPracticeQns.Range("B10").NumberFormat = "@"

# Switch the active sheet from PracticeQns to TryCode, with a new selection
[void]$wsTryCode.Activate()
[void]$wsTryCode.Range("L20").Select()
